# Update cryptocurrency price/volume/date/hour data to the Feb 12 2023 snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "309.69"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "1.05%"
$c = $ws.Range("F2")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "41.09"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "0.78%"
$c = $ws.Range("F3")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.149"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "0.17%"
$c = $ws.Range("F4")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.07687"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "1.04%"
$c = $ws.Range("F5")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.618"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "0.54%"
$c = $ws.Range("F6")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.9226"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "1.85%"
$c = $ws.Range("F7")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "2.481"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "0.61%"
$c = $ws.Range("F8")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1160"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "11.97%"
$c = $ws.Range("F9")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "5.20%"
$c = $ws.Range("F10")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.09167"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "-0.11%"
$c = $ws.Range("F11")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.04293"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "1.66%"
$c = $ws.Range("F12")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-0.71%"
$c = $ws.Range("F13")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001257"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "1.00%"
$c = $ws.Range("F14")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.005791"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "-1.10%"
$c = $ws.Range("F15")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.343"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "-0.25%"
$c = $ws.Range("F16")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.328"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "1.64%"
$c = $ws.Range("F17")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.3335"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "1.88%"
$c = $ws.Range("F18")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.030"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "7.06%"
$c = $ws.Range("F19")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.1400"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "2.52%"
$c = $ws.Range("F20")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.2913"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "7.02%"
$c = $ws.Range("F21")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.04055"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "-2.77%"
$c = $ws.Range("F22")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.001263"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "3.15%"
$c = $ws.Range("F23")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.004123"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "1.30%"
$c = $ws.Range("F24")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0001272"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "-2.26%"
$c = $ws.Range("F25")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "24.54%"
$c = $ws.Range("F26")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F27")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F28")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F29")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F30")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F31")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F32")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F33")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F34")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F35")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F36")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("F37")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02445"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "3.00%"
$c = $ws.Range("F38")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05280"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "2.60%"
$c = $ws.Range("F39")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.007844"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "0.81%"
$c = $ws.Range("F40")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1314"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "1.42%"
$c = $ws.Range("F41")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.006808"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "-2.14%"
$c = $ws.Range("F42")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.001902"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "-0.93%"
$c = $ws.Range("F43")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.007430"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-14.95%"
$c = $ws.Range("F44")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.3371"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "1.42%"
$c = $ws.Range("F45")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00006805"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "6.96%"
$c = $ws.Range("F46")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.00000000751"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "-0.04%"
$c = $ws.Range("F47")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.1699"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "2,109.17%"
$c = $ws.Range("F48")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.004100"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "-6.97%"
$c = $ws.Range("F49")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.00002102"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.04%"
$c = $ws.Range("F50")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = "0"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0002002"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "-0.04%"
$c = $ws.Range("F51")
$c.NumberFormat = "@"
$c.Value = "12-2-2023"
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = "0"
